# Updates the "cryptos" price/volume table (Sheet1) to reflect the latest
# scrape of coinranking.com data. For each affected row we overwrite only
# the columns that changed: Coin (B), Link (C), Price (D) and Volume(1h) (E).
#
# Column D values are plain text that often look like numbers containing a
# single "." (e.g. "7.90", "2.01") or look like grouped numbers with more
# than one "." (e.g. "51.823.86"). Assigning such a string straight to
# .Value would let Excel auto-convert it into a real number (dropping
# trailing zeros, truncating to a double, etc.), so for column D we force
# the cell to Text format ("@") before writing the value and then restore
# the default "Normal" style afterwards so the cell's formatting is left
# exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; D = '51.823.86'; E = '  -0.46%  ' }
    @{ Row = 3; D = '2.928.25'; E = '  +1.46%  ' }
    @{ Row = 4; E = '  -0.11%  ' }
    @{ Row = 5; D = '355.03'; E = '  +1.13%  ' }
    @{ Row = 6; D = '110.84'; E = '  -0.53%  ' }
    @{ Row = 7; D = '0.568'; E = '  +2.11%  ' }
    @{ Row = 8; E = '  -0.01%  ' }
    @{ Row = 9; D = '0.629'; E = '  +1.21%  ' }
    @{ Row = 10; D = '39.33'; E = '  -1.24%  ' }
    @{ Row = 11; D = '0.0882'; E = '  +3.25%  ' }
    @{ Row = 12; D = '0.136'; E = '  +0.57%  ' }
    @{ Row = 13; E = '  -1.02%  ' }
    @{ Row = 14; D = '7.90'; E = '  +1.66%  ' }
    @{ Row = 15; D = '3.385.08'; E = '  +1.17%  ' }
    @{ Row = 16; D = '2.945.72'; E = '  +1.48%  ' }
    @{ Row = 17; D = '0.985'; E = '  -1.53%  ' }
    @{ Row = 18; D = '51.813.00'; E = '  -0.49%  ' }
    @{ Row = 19; E = '  -1.59%  ' }
    @{ Row = 20; D = '7.55'; E = '  -1.94%  ' }
    @{ Row = 21; D = '14.04'; E = '  -2.90%  ' }
    @{ Row = 22; E = '  +0.43%  ' }
    @{ Row = 23; D = '70.88'; E = '  +0.17%  ' }
    @{ Row = 24; D = '270.70'; E = '  +0.46%  ' }
    @{ Row = 25; D = '2.81'; E = '  +0.77%  ' }
    @{ Row = 26; D = '0.182'; E = '  +11.50%  ' }
    @{ Row = 27; D = '27.15'; E = '  +2.69%  ' }
    @{ Row = 28; E = '  +0.22%  ' }
    @{ Row = 29; D = '7.32'; E = '  +13.95%  ' }
    @{ Row = 30; D = '0.107'; E = '  +12.96%  ' }
    @{ Row = 31; B = 'InjectiveProtocol'; C = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D = '39.05'; E = '  +1.78%  ' }
    @{ Row = 32; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '10.57'; E = '  +0.81%  ' }
    @{ Row = 33; D = '6.06'; E = '  -1.46%  ' }
    @{ Row = 34; D = '52.26'; E = '  -1.16%  ' }
    @{ Row = 35; D = '0.0442'; E = '  -3.67%  ' }
    @{ Row = 36; E = '  -0.02%  ' }
    @{ Row = 37; E = '  -14.77%  ' }
    @{ Row = 38; D = '3.24'; E = '  -1.93%  ' }
    @{ Row = 39; D = '18.46'; E = '  -0.79%  ' }
    @{ Row = 40; B = 'ARBITRUM'; C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D = '2.01'; E = '  -1.37%  ' }
    @{ Row = 41; B = 'Stacks'; C = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D = '2.77'; E = '  +4.40%  ' }
    @{ Row = 42; E = '  +2.72%  ' }
    @{ Row = 43; D = '23.22'; E = '  +1.97%  ' }
    @{ Row = 44; D = '120.21'; E = '  -1.40%  ' }
    @{ Row = 45; E = '  -1.71%  ' }
    @{ Row = 46; E = '  +0.46%  ' }
    @{ Row = 47; D = '3.46'; E = '  -3.16%  ' }
    @{ Row = 48; D = '2.136.44'; E = '  -3.00%  ' }
    @{ Row = 49; E = '  -6.36%  ' }
    @{ Row = 50; D = '0.0335'; E = '  +4.49%  ' }
    @{ Row = 51; B = 'FraxShare'; C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D = '9.11'; E = '  +0.36%  ' }
)

foreach ($item in $rows) {
    $r = $item.Row

    if ($item.ContainsKey('B')) {
        $ws.Range("B$r").Value = $item.B
    }
    if ($item.ContainsKey('C')) {
        $ws.Range("C$r").Value = $item.C
    }
    if ($item.ContainsKey('D')) {
        # Force text storage so numeric-looking strings (e.g. "7.90",
        # "51.823.86") are not reinterpreted as numbers/dates by Excel.
        $cell = $ws.Range("D$r")
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
        $cell.Style = "Normal"
    }
    if ($item.ContainsKey('E')) {
        $ws.Range("E$r").Value = $item.E
    }
}
